$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ws, $ref, $val) {
    $rng = $ws.Range($ref)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextCell $ws "D2" "258.48"
Set-TextCell $ws "E2" "0.49%"
Set-TextCell $ws "G2" "19"
Set-TextCell $ws "D3" "26.90"
Set-TextCell $ws "E3" "-1.74%"
Set-TextCell $ws "G3" "19"
Set-TextCell $ws "D4" "4.653"
Set-TextCell $ws "E4" "1.36%"
Set-TextCell $ws "G4" "19"
Set-TextCell $ws "D5" "0.05997"
Set-TextCell $ws "E5" "1.84%"
Set-TextCell $ws "G5" "19"
Set-TextCell $ws "D6" "6.667"
Set-TextCell $ws "E6" "0.60%"
Set-TextCell $ws "G6" "19"
Set-TextCell $ws "D7" "0.8584"
Set-TextCell $ws "E7" "0.05%"
Set-TextCell $ws "G7" "19"
Set-TextCell $ws "D8" "0.9218"
Set-TextCell $ws "E8" "-0.31%"
Set-TextCell $ws "G8" "19"
Set-TextCell $ws "D9" "0.1393"
Set-TextCell $ws "E9" "-1.21%"
Set-TextCell $ws "G9" "19"
Set-TextCell $ws "D10" "0.04524"
Set-TextCell $ws "E10" "26.83%"
Set-TextCell $ws "G10" "19"
Set-TextCell $ws "D11" "0.07028"
Set-TextCell $ws "E11" "-0.83%"
Set-TextCell $ws "G11" "19"
Set-TextCell $ws "D12" "0.03116"
Set-TextCell $ws "E12" "-2.99%"
Set-TextCell $ws "G12" "19"
Set-TextCell $ws "D13" "0.09132"
Set-TextCell $ws "E13" "-0.50%"
Set-TextCell $ws "G13" "19"
Set-TextCell $ws "D14" "0.001527"
Set-TextCell $ws "E14" "-1.11%"
Set-TextCell $ws "G14" "19"
Set-TextCell $ws "D15" "0.0006049"
Set-TextCell $ws "E15" "-94.21%"
Set-TextCell $ws "G15" "19"
Set-TextCell $ws "D16" "0.006152"
Set-TextCell $ws "E16" "2.40%"
Set-TextCell $ws "G16" "19"
Set-TextCell $ws "E17" "-1.68%"
Set-TextCell $ws "G17" "19"
Set-TextCell $ws "D18" "3.160"
Set-TextCell $ws "E18" "-1.30%"
Set-TextCell $ws "G18" "19"
Set-TextCell $ws "E19" "-0.82%"
Set-TextCell $ws "G19" "19"
Set-TextCell $ws "E20" "0.10%"
Set-TextCell $ws "G20" "19"
Set-TextCell $ws "D21" "0.1288"
Set-TextCell $ws "E21" "0.85%"
Set-TextCell $ws "G21" "19"
Set-TextCell $ws "D22" "4.153"
Set-TextCell $ws "E22" "7.84%"
Set-TextCell $ws "G22" "19"
Set-TextCell $ws "D23" "0.04238"
Set-TextCell $ws "E23" "0.89%"
Set-TextCell $ws "G23" "19"
Set-TextCell $ws "D24" "0.001216"
Set-TextCell $ws "E24" "-0.22%"
Set-TextCell $ws "G24" "19"
Set-TextCell $ws "D25" "0.004033"
Set-TextCell $ws "E25" "-6.18%"
Set-TextCell $ws "G25" "19"
Set-TextCell $ws "D26" "0.0001200"
Set-TextCell $ws "E26" "0.05%"
Set-TextCell $ws "G26" "19"
Set-TextCell $ws "D27" "0.0001714"
Set-TextCell $ws "E27" "13.55%"
Set-TextCell $ws "G27" "19"
Set-TextCell $ws "G28" "19"
Set-TextCell $ws "G29" "19"
Set-TextCell $ws "G30" "19"
Set-TextCell $ws "G31" "19"
Set-TextCell $ws "G32" "19"
Set-TextCell $ws "G33" "19"
Set-TextCell $ws "G34" "19"
Set-TextCell $ws "G35" "19"
Set-TextCell $ws "G36" "19"
Set-TextCell $ws "G37" "19"
Set-TextCell $ws "G38" "19"
Set-TextCell $ws "G39" "19"
Set-TextCell $ws "D40" "0.03840"
Set-TextCell $ws "E40" "0.20%"
Set-TextCell $ws "G40" "19"
Set-TextCell $ws "D41" "0.1115"
Set-TextCell $ws "E41" "1.14%"
Set-TextCell $ws "G41" "19"
Set-TextCell $ws "D42" "0.003850"
Set-TextCell $ws "E42" "-38.45%"
Set-TextCell $ws "G42" "19"
Set-TextCell $ws "D43" "0.002420"
Set-TextCell $ws "E43" "10.05%"
Set-TextCell $ws "G43" "19"
Set-TextCell $ws "D44" "0.01527"
Set-TextCell $ws "E44" "29.31%"
Set-TextCell $ws "G44" "19"
Set-TextCell $ws "D45" "0.00005092"
Set-TextCell $ws "E45" "-6.74%"
Set-TextCell $ws "G45" "19"
Set-TextCell $ws "E46" "0.04%"
Set-TextCell $ws "G46" "19"
Set-TextCell $ws "E47" "-20.66%"
Set-TextCell $ws "G47" "19"
Set-TextCell $ws "E48" "1.26%"
Set-TextCell $ws "G48" "19"
Set-TextCell $ws "D49" "0.00002100"
Set-TextCell $ws "E49" "0.04%"
Set-TextCell $ws "G49" "19"
Set-TextCell $ws "D50" "0.0002000"
Set-TextCell $ws "E50" "0.04%"
Set-TextCell $ws "G50" "19"
Set-TextCell $ws "G51" "19"
